# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 71 (pushing the existing
# rows 71-98 down to 72-99); the new row carries a new "Red Beaut" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71; everything below (old rows 71-98) shifts
# down by one (to 72-99), carrying its formatting (incl. the date style).
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's data.
$ws.Range("A71").Value = 4
$ws.Range("B71").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value = "Los Lagos"
$ws.Range("D71").Value = 44553
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100103
$ws.Range("H71").Value = "Frutos de hueso (carozo)"
$ws.Range("I71").Value = 100103002
$ws.Range("J71").Value = "Ciruela"
$ws.Range("K71").Value = "Red Beaut"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 700
$ws.Range("N71").Value = 19000
$ws.Range("O71").Value = 20000
$ws.Range("P71").Value = 19500
$ws.Range("Q71").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R71").Value = "Región Metropolitana"
$ws.Range("S71").Value = 1300
$ws.Range("T71").Value = 15
